$wb = $excel.ActiveWorkbook

$replacements = @(
    @("Tom's Baked Beans", "Toms Baked Beans"),
    @("Haley's David", "Haleys David"),
    @("Sun's Teriyaki", "Suns Teriyaki"),
    @("BYU-Idaho", "BYU Idaho"),
    @("Everybody's a Lobo", "Everybodys a Lobo"),
    @("Megan's Megans", "Megans Megans"),
    @("Dad's friend Mike", "Dads friend Mike"),
    @("Trick 'r Treat", "Trick r Treat"),
    @("Child's Play", "Childs Play"),
    @("You're Next", "Youre Next"),
    @("Jennifer's Body", "Jennifers Body"),
    @("Well I don't approve", "Well I dont approve"),
    @("I don't like sand.", "I dont like sand."),
    @("It's course and rough and irritating", "Its course and rough and irritating"),
    @("Bad bombin'", "Bad bombin"),
    @("That's why I'm here", "Thats why Im here"),
    @("It's an older code sir, but it checks out", "Its an older code sir, but it checks out"),
    @("Eznog's Eggnog", "Eznogs Eggnog"),
    @("Beau's PBS show", "Beaus PBS show"),
    @("Miner's Treat", "Miners Treat"),
    @("Lil' Sebastian", "Lil Sebastian"),
    @("Tom's Bistro", "Toms Bistro"),
    @("JJ's Diner", "JJs Diner"),
    @("There's always money in the banana stand", "Theres always money in the banana stand"),
    @("He's going to be all right", "Hes going to be all right"),
    @("Alice's Adventures in Wonderland", "Alices Adventures in Wonderland"),
    @("Gulliver's Travels", "Gullivers Travels"),
    @("The Pilgrim's Progress", "The Pilgrims Progress"),
    @("Uncle Tom's Cabin", "Uncle Toms Cabin"),
    @("Tess of the d'Urbervilles", "Tess of the dUrbervilles"),
    @("Harry Potter and the Sorcerer's Stone", "Harry Potter and the Sorcerers Stone"),
    @("Lady Chatterley's Lover", "Lady Chatterleys Lover"),
    @("A Connecticut Yankee in King Arthur's Court", "A Connecticut Yankee in King Arthurs Court"),
    @("Charlotte's Web", "Charlottes Web"),
    @("The Magician's Nephew", "The Magicians Nephew"),
    @("One Flew Over the Cuckoo's Nest", "One Flew Over the Cuckoos Nest"),
    @("The Handmaid's Tale", "The Handmaids Tale"),
    @("All the King's Men", "All the Kings Men"),
    @("Are You There God? It's Me, Margaret", "Are You There God? Its Me, Margaret"),
    @("Swann's Way", "Swanns Way"),
    @("Ender's Game", "Enders Game"),
    @("The Hitchhiker's Guide to the Galaxy", "The Hitchhikers Guide to the Galaxy"),
    @("Le Morte D'Arthur", "Le Morte DArthur"),
    @("The Girl Who Kicked the Hornet's Nest", "The Girl Who Kicked the Hornets Nest"),
    @("My Sister's Keeper", "My Sisters Keeper"),
    @("Schindler's List", "Schindlers List"),
    @("The Time Traveler's Wife", "The Time Travelers Wife"),
    @("King Solomon's Mines", "King Solomons Mines"),
    @("The No. 1 Ladies' Detective Agency", "The No. 1 Ladies Detective Agency"),
    @("Miss Peregrine's Home for Peculiar Children", "Miss Peregrines Home for Peculiar Children"),
    @("Sophie's Choice", "Sophies Choice"),
    @("Breakfast at Tiffany's", "Breakfast at Tiffanys"),
    @("The Razor's Edge", "The Razors Edge"),
    @("The Memory Keeper's Daughter", "The Memory Keepers Daughter"),
    @("Alex Cross's Trial", "Alex Crosss Trial"),
    @("Cat's Cradle", "Cats Cradle"),
    @("Madeleine L'Engle", "Madeleine LEngle"),
    @("Robert C. O'Brien", "Robert C. OBrien"),
    @("Tim O'Brien", "Tim OBrien"),
    @("It's Always Sunny in Philadelphia", "Its Always Sunny in Philadelphia"),
    @("Grey's Anatomy", "Greys Anatomy"),
    @("Handmaid's Tale", "Handmaids Tale")
)

foreach ($ws in $wb.Worksheets) {
    foreach ($pair in $replacements) {
        $find = $pair[0]
        $replace = $pair[1]
        [void]$ws.Cells.Replace($find, $replace, 1, 1, $false, $false, $false, $false)
    }
}
